$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1 - Arrival_rates")

# Replace the placeholder "X" values with "0" (kept as text, matching the
# cells' Text number format) across the used data range.
$ws.Range("C2").Value = "0"
$ws.Range("C3").Value = "0"
$ws.Range("D3").Value = "0"
$ws.Range("E3").Value = "0"
$ws.Range("D4").Value = "0"
$ws.Range("E4").Value = "0"

# Update the active selection to N7, as recorded in the saved view state.
$ws.Range("N7").Select()
